## Updated Todo and README with latest build
##
## "Immediate Checklist" sheet: several completed/duplicate items are
## removed. The "Ammo starts on 0/0..." item (previously row 4) is the
## only one that survives, and it is moved up into row 2 (replacing the
## now-finished "Add sound effects..." item). The remaining rows that
## held now-removed items are cleared out. Row heights that were only
## custom because of the (now gone) two-line text are reset back to the
## sheet's default via AutoFit. The shared strings for the four removed
## items are dropped automatically on save once nothing references them.
##
## "Features" sheet: no content changed, just where the cursor was left.

$wb = $excel.ActiveWorkbook

$checklist = $wb.Worksheets.Item("Immediate Checklist")

# Row 2 used to read "Add sound effects for bullet chambering and enemies
# attacking" - that task is done, so the still-open "Ammo starts on 0/0..."
# item (which used to live on row 4) takes its place.
$checklist.Range("A2").Value = "Ammo starts on 0/0 for some reason in builds."

# Row 3 ("Play a different sound in the shop...") is finished - clear it.
$checklist.Range("A3").ClearContents()
$checklist.Rows(3).AutoFit()

# Row 4 used to hold "Ammo starts on 0/0..." which has now been moved up
# into row 2, so this row is emptied out.
$checklist.Range("A4").ClearContents()
$checklist.Rows(4).AutoFit()

# Row 5 ("Lock the enemy health bars to only rotate on y-axis.") is done.
$checklist.Range("A5").ClearContents()
$checklist.Rows(5).AutoFit()

# Row 6 ("Dealing 100 damange with sniper rifle...") is done.
$checklist.Range("A6").ClearContents()
$checklist.Rows(6).AutoFit()

# Features sheet: just move the selection, no data changes.
$features = $wb.Worksheets.Item("Features")
$features.Range("A3").Select()

# Leave the cursor where the last edit on the checklist sheet happened,
# and leave that sheet as the active tab (matches the saved workbook).
$checklist.Range("A14").Select()
